# Apply "Dish - Verified Code" edits
$wb = $excel.ActiveWorkbook

$wsCode = $wb.Worksheets.Item("Code Parts")
$wsMethods = $wb.Worksheets.Item("Methods Required")

# --- Sheet "Code Parts" (sheet1) ---
# Row 5: C5:F5 -> "x"
$wsCode.Range("C5:F5").Value = "x"
# Row 9: E9:F9 -> "x"
$wsCode.Range("E9:F9").Value = "x"

# Column F width + selection/zoom
$wsCode.Columns.Item(6).ColumnWidth = 18.140625
$wsCode.Range("E13").Select()
$wsCode.Application.ActiveWindow.Zoom = 115

# --- Sheet "Methods Required" (sheet2) ---
$wsMethods.Range("C4").Value = "x"
$wsMethods.Range("C5:D5").Value = "x"
$wsMethods.Range("C6:D6").Value = "x"
$wsMethods.Range("C7:D7").Value = "x"
$wsMethods.Range("C8:D8").Value = "x"
$wsMethods.Range("C22:D22").Value = "x"
$wsMethods.Range("C23").Value = "x"

$wsMethods.Activate()
$wsMethods.Range("G7").Select()
$wsMethods.Application.ActiveWindow.Zoom = 100
